# fix: validasi input, lokasi dosen management
#
# Template perubahan: kolom "Kategori" pada baris contoh (F2) diubah dari
# "Nasional" menjadi "Akademik".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Kategori contoh lomba: "Nasional" -> "Akademik"
$ws.Range("F2").Value = "Akademik"

# Touch the header/footer settings so the sheet carries an explicit (empty)
# headerFooter section, matching the regenerated template's page setup block.
$ps = $ws.PageSetup
$ps.LeftHeader = ""
$ps.CenterHeader = ""
$ps.RightHeader = ""
$ps.LeftFooter = ""
$ps.CenterFooter = ""
$ps.RightFooter = ""
